$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# like "0.999" are not auto-converted to numbers, matching the
# original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '73.091.57'
$ws.Range("E2").Value = '  +2.95%  '
$ws.Range("D3").Value = '3.993.26'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '595.08'
$ws.Range("E5").Value = '  +10.79%  '
$ws.Range("D6").Value = '164.66'
$ws.Range("E6").Value = '  +11.27%  '
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.751'
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("D11").Value = '54.81'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("E13").Value = '  +3.46%  '
$ws.Range("D14").Value = '4.632.43'
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").Value = '3.998.50'
$ws.Range("E16").Value = '  +8.68%  '
$ws.Range("D17").Value = '14.16'
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").Value = '20.47'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").Value = '72.790.20'
$ws.Range("E20").Value = '  +2.65%  '
$ws.Range("D21").Value = '438.41'
$ws.Range("E21").Value = '  +4.05%  '
$ws.Range("D22").Value = '4.75'
$ws.Range("E22").Value = '  +12.23%  '
$ws.Range("D23").Value = '96.67'
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("E24").Value = '  -4.27%  '
$ws.Range("D25").Value = '14.38'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '4.36'
$ws.Range("E26").Value = '  +14.50%  '
$ws.Range("D27").Value = '11.33'
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("D29").Value = '10.41'
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").Value = '36.44'
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").Value = '13.95'
$ws.Range("E32").Value = '  +4.65%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '48.42'
$ws.Range("E34").Value = '  -5.49%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '670.01'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").Value = '70.73'
$ws.Range("E36").Value = '  +8.09%  '
$ws.Range("D37").Value = '0.0₃0903'
$ws.Range("E37").Value = '  +10.59%  '
$ws.Range("D38").Value = '0.439'
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = '3.36'
$ws.Range("E41").Value = '  +5.47%  '
$ws.Range("D42").Value = '0.146'
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '0.0491'
$ws.Range("E44").Value = '  +2.08%  '
$ws.Range("D45").Value = '10.69'
$ws.Range("E45").Value = '  +6.25%  '
$ws.Range("D46").Value = '0.150'
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").Value = '3.41'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").Value = '2.925.24'
$ws.Range("E49").Value = '  +11.35%  '
$ws.Range("E50").Value = '  +2.23%  '
$ws.Range("D51").Value = '3.41'
$ws.Range("E51").Value = '  +4.67%  '

# Restore default (unstyled) formatting on column D now that the
# text values are safely written, so no stray number-format is left
# applied to the cells.
$ws.Range("D2:D51").ClearFormats()
